$p = $ppt.ActivePresentation

# Append one new slide at the end of the deck, using the "Blank" layout
# (ppLayoutBlank = 12), matching every other slide that was freshly added
# to this deck (slides 2-5 all use the Blank layout already).
#
# We duplicate the current last slide (it already uses the Blank layout)
# rather than calling Slides.Add directly, because duplicating carries
# over the full <p:cSld>/<p:grpSpPr>/<p:extLst>/<p:clrMapOvr> scaffolding
# that PowerPoint itself stamps onto a slide, instead of the bare-bones
# tree a plain Add would synthesize. We then strip out the shapes that
# came along with the duplicate so the appended slide ends up completely
# blank, and make sure its layout is still Blank.
$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastSlide.Duplicate() | Out-Null
$newSlide = $p.Slides.Item($p.Slides.Count)

$newSlide.Layout = 12

while ($newSlide.Shapes.Count -gt 0) {
    $newSlide.Shapes.Item(1).Delete()
}
